$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.133572101593018
$ws.Range("B1").Value = 2.328168630599976
$ws.Range("C1").Value = 1.936259865760803
$ws.Range("D1").Value = 1.84643828868866
$ws.Range("E1").Value = 1.653774380683899
